# "pulls most recent US EPS commit"
#
# Updates the Capacity Supply Curve workbook to match the latest upstream
# data: the onshore-wind "share of existing capacity built this year" row
# drops from 0.3 to 0.2 for every modeled year, and the workbook is left
# with the data sheet active/selected (reflecting the author's last view)
# instead of the About sheet.

$wb = $excel.ActiveWorkbook

# --- Data change -------------------------------------------------------
# Sheet "CSC-CSCSoCECBiaSY" ("Share of Cost Effective Capacity Built in a
# Single Year"), row 7 = "onshore wind es": 0.3 -> 0.2 for every year
# column (B:AE).
$wsShare = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$wsShare.Range("B7:AE7").Value = 0.2

# --- Active sheet / selection -------------------------------------------
# Move off of "About" (no longer the selected tab) and onto the data
# sheet just edited, selecting the row that was changed.
$wsShare.Activate()
$wsShare.Range("B7:AE7").Select()
